$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 45096
$ws.Cells.Item(2, 13).Value = 50
$ws.Cells.Item(2, 14).Value = 23000
$ws.Cells.Item(2, 15).Value = 23000
$ws.Cells.Item(2, 16).Value = 23000
$ws.Cells.Item(2, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(2, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(2, 19).Value = 1278
$ws.Cells.Item(2, 20).Value = 18

$ws.Cells.Item(3, 4).Value = 45014
$ws.Cells.Item(3, 13).Value = 30
$ws.Cells.Item(3, 14).Value = 18000
$ws.Cells.Item(3, 15).Value = 18000
$ws.Cells.Item(3, 16).Value = 18000
$ws.Cells.Item(3, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(3, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(3, 19).Value = 1000
$ws.Cells.Item(3, 20).Value = 18

$ws.Cells.Item(4, 4).Value = 45089
$ws.Cells.Item(4, 13).Value = 60
$ws.Cells.Item(4, 14).Value = 22000
$ws.Cells.Item(4, 15).Value = 23000
$ws.Cells.Item(4, 16).Value = 22500
$ws.Cells.Item(4, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(4, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(4, 19).Value = 1250
$ws.Cells.Item(4, 20).Value = 18

$ws.Cells.Item(5, 4).Value = 45062
$ws.Cells.Item(5, 13).Value = 90
$ws.Cells.Item(5, 14).Value = 13000
$ws.Cells.Item(5, 15).Value = 14000
$ws.Cells.Item(5, 16).Value = 13444
$ws.Cells.Item(5, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(5, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(5, 19).Value = 747
$ws.Cells.Item(5, 20).Value = 18

$ws.Cells.Item(6, 4).Value = 45028
$ws.Cells.Item(6, 13).Value = 50
$ws.Cells.Item(6, 14).Value = 18000
$ws.Cells.Item(6, 15).Value = 18000
$ws.Cells.Item(6, 16).Value = 18000
$ws.Cells.Item(6, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(6, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(6, 19).Value = 1000
$ws.Cells.Item(6, 20).Value = 18

$ws.Cells.Item(7, 4).Value = 45043
$ws.Cells.Item(7, 13).Value = 60
$ws.Cells.Item(7, 14).Value = 15000
$ws.Cells.Item(7, 15).Value = 15000
$ws.Cells.Item(7, 16).Value = 15000
$ws.Cells.Item(7, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(7, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(7, 19).Value = 833
$ws.Cells.Item(7, 20).Value = 18

$ws.Cells.Item(8, 4).Value = 45037
$ws.Cells.Item(8, 13).Value = 60
$ws.Cells.Item(8, 14).Value = 16000
$ws.Cells.Item(8, 15).Value = 16000
$ws.Cells.Item(8, 16).Value = 16000
$ws.Cells.Item(8, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(8, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(8, 19).Value = 889
$ws.Cells.Item(8, 20).Value = 18

$ws.Cells.Item(9, 4).Value = 44999
$ws.Cells.Item(9, 13).Value = 60
$ws.Cells.Item(9, 14).Value = 17000
$ws.Cells.Item(9, 15).Value = 18000
$ws.Cells.Item(9, 16).Value = 17500
$ws.Cells.Item(9, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(9, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(9, 19).Value = 972
$ws.Cells.Item(9, 20).Value = 18

$ws.Cells.Item(10, 4).Value = 45091
$ws.Cells.Item(10, 13).Value = 50
$ws.Cells.Item(10, 14).Value = 22000
$ws.Cells.Item(10, 15).Value = 22000
$ws.Cells.Item(10, 16).Value = 22000
$ws.Cells.Item(10, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(10, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(10, 19).Value = 1222
$ws.Cells.Item(10, 20).Value = 18

$ws.Cells.Item(11, 4).Value = 45041
$ws.Cells.Item(11, 13).Value = 60
$ws.Cells.Item(11, 14).Value = 15000
$ws.Cells.Item(11, 15).Value = 15000
$ws.Cells.Item(11, 16).Value = 15000
$ws.Cells.Item(11, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(11, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(11, 19).Value = 833
$ws.Cells.Item(11, 20).Value = 18

$ws.Cells.Item(12, 4).Value = 45033
$ws.Cells.Item(12, 13).Value = 60
$ws.Cells.Item(12, 14).Value = 15000
$ws.Cells.Item(12, 15).Value = 16000
$ws.Cells.Item(12, 16).Value = 15500
$ws.Cells.Item(12, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(12, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(12, 19).Value = 861
$ws.Cells.Item(12, 20).Value = 18

$ws.Cells.Item(13, 4).Value = 45020
$ws.Cells.Item(13, 13).Value = 50
$ws.Cells.Item(13, 14).Value = 15000
$ws.Cells.Item(13, 15).Value = 15000
$ws.Cells.Item(13, 16).Value = 15000
$ws.Cells.Item(13, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(13, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(13, 19).Value = 938
$ws.Cells.Item(13, 20).Value = 16

$ws.Cells.Item(14, 4).Value = 45021
$ws.Cells.Item(14, 13).Value = 60
$ws.Cells.Item(14, 14).Value = 15000
$ws.Cells.Item(14, 15).Value = 16000
$ws.Cells.Item(14, 16).Value = 15500
$ws.Cells.Item(14, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(14, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(14, 19).Value = 861
$ws.Cells.Item(14, 20).Value = 18

$ws.Cells.Item(15, 4).Value = 45044
$ws.Cells.Item(15, 13).Value = 60
$ws.Cells.Item(15, 14).Value = 15000
$ws.Cells.Item(15, 15).Value = 15000
$ws.Cells.Item(15, 16).Value = 15000
$ws.Cells.Item(15, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(15, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(15, 19).Value = 833
$ws.Cells.Item(15, 20).Value = 18

$ws.Cells.Item(16, 4).Value = 45099
$ws.Cells.Item(16, 13).Value = 40
$ws.Cells.Item(16, 14).Value = 22000
$ws.Cells.Item(16, 15).Value = 22000
$ws.Cells.Item(16, 16).Value = 22000
$ws.Cells.Item(16, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(16, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(16, 19).Value = 1222
$ws.Cells.Item(16, 20).Value = 18

$ws.Cells.Item(17, 4).Value = 45001
$ws.Cells.Item(17, 13).Value = 60
$ws.Cells.Item(17, 14).Value = 17000
$ws.Cells.Item(17, 15).Value = 18000
$ws.Cells.Item(17, 16).Value = 17500
$ws.Cells.Item(17, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(17, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(17, 19).Value = 972
$ws.Cells.Item(17, 20).Value = 18

$ws.Cells.Item(18, 4).Value = 45036
$ws.Cells.Item(18, 13).Value = 60
$ws.Cells.Item(18, 14).Value = 15000
$ws.Cells.Item(18, 15).Value = 16000
$ws.Cells.Item(18, 16).Value = 15500
$ws.Cells.Item(18, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(18, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(18, 19).Value = 861
$ws.Cells.Item(18, 20).Value = 18

$ws.Cells.Item(19, 4).Value = 45050
$ws.Cells.Item(19, 13).Value = 40
$ws.Cells.Item(19, 14).Value = 14000
$ws.Cells.Item(19, 15).Value = 14000
$ws.Cells.Item(19, 16).Value = 14000
$ws.Cells.Item(19, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(19, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(19, 19).Value = 778
$ws.Cells.Item(19, 20).Value = 18

$ws.Cells.Item(20, 4).Value = 45049
$ws.Cells.Item(20, 13).Value = 80
$ws.Cells.Item(20, 14).Value = 15000
$ws.Cells.Item(20, 15).Value = 15000
$ws.Cells.Item(20, 16).Value = 15000
$ws.Cells.Item(20, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(20, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(20, 19).Value = 833
$ws.Cells.Item(20, 20).Value = 18

$ws.Cells.Item(21, 4).Value = 45030
$ws.Cells.Item(21, 13).Value = 40
$ws.Cells.Item(21, 14).Value = 18000
$ws.Cells.Item(21, 15).Value = 18000
$ws.Cells.Item(21, 16).Value = 18000
$ws.Cells.Item(21, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(21, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(21, 19).Value = 1000
$ws.Cells.Item(21, 20).Value = 18

$ws.Cells.Item(22, 4).Value = 45002
$ws.Cells.Item(22, 13).Value = 30
$ws.Cells.Item(22, 14).Value = 18000
$ws.Cells.Item(22, 15).Value = 18000
$ws.Cells.Item(22, 16).Value = 18000
$ws.Cells.Item(22, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(22, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(22, 19).Value = 1000
$ws.Cells.Item(22, 20).Value = 18
